$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the Jan 28 2024 GitHub Actions refresh.
# Column D ("Price") values that look numeric must be forced to remain plain text
# (matching the source inlineStr cells, which preserve formatting like trailing zeros,
# e.g. "156.30" rather than being re-interpreted as the number 156.3).

$ws.Range('D2').Value = '42.694.12'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '2.300.99'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '156.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +15,508.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '96.72'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.75%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.499'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.83%  '
$ws.Range('E11').Value = '  +8.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0812'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E13').Value = '  -1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('D15').Value = '2.657.09'
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '2.310.90'
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.798'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.17%  '
$ws.Range('D19').Value = '42.557.09'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.92%  '
$ws.Range('D21').Value = '0.0₃0923'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '244.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.54%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.75'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('E31').Value = '  +0.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('E33').Value = '  +4.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0760'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.36%  '
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.110'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.07%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.42'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.117'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.68%  '
$ws.Range('D44').Value = '2.022.22'
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range('E45').Value = '  +10.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0287'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.55%  '
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.03'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.55'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.75'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.17%  '
